# Generate Report for Handoff
# Updates the localization-status workbook after a new handoff for b.md:
#   - Overview sheet: b.md row now shows "Ready for handoff" status for both
#     locales, and a refreshed "Latest HO Xliff Generate Date".
#   - zh-cn / de-de sheets: b.md row gets a new Status, Content Duplicate flips
#     to False, a new Latest Handoff File / Datetime, and an Error Detail
#     message about the handback file being stale. Column P (Error Detail) is
#     widened to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$newDate = "2016-08-31 00:40:30"
$xlPasteValues = -4163

# ---------------------------------------------------------------------------
# Overview sheet (row 3 = b.md)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("G3").Value = $newDate

# ---------------------------------------------------------------------------
# zh-cn sheet (row 3 = b.md)
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus

# "Content Duplicate" is a text column holding the literal words True/False
# (see F2, which already stores text "False"). Assigning the string "False"
# directly to .Value gets auto-detected as a native boolean, so instead copy
# the already-text-typed neighbour cell and paste-special just the value,
# which preserves the text type.
$zhcn.Range("F2").Copy()
$zhcn.Range("F3").PasteSpecial($xlPasteValues)
$excel.CutCopyMode = $false

$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-31 00:40:26"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/64e1ca09ac2bd913f8b030dac40b2055571fed10/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4e4c9045a5fa00afc84c4abc70b57293545fe4b7/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.1428571428571

# ---------------------------------------------------------------------------
# de-de sheet (row 3 = b.md)
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus

$dede.Range("F2").Copy()
$dede.Range("F3").PasteSpecial($xlPasteValues)
$excel.CutCopyMode = $false

$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = $newDate
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/64e1ca09ac2bd913f8b030dac40b2055571fed10/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4e4c9045a5fa00afc84c4abc70b57293545fe4b7/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.1428571428571

Write-Output "Report for handoff generated"
